# Update the "想去人数" (want-to-go count) column F on each sheet to the
# refreshed snapshot values captured in the new data pull.
#
# Sheet "展览"   (sheet1 / index 1)
# Sheet "演出"   (sheet2 / index 2)
# Sheet "本地生活" (sheet3 / index 3)
# Sheet "全部类型" (sheet4 / index 4)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 37723
$ws1.Range("F9").Value = 852
$ws1.Range("F22").Value = 845
$ws1.Range("F23").Value = 2552
$ws1.Range("F24").Value = 1028
$ws1.Range("F29").Value = 793
$ws1.Range("F31").Value = 1168

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 334

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 641

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 641
$ws4.Range("F3").Value = 37723
$ws4.Range("F12").Value = 334
$ws4.Range("F15").Value = 852
$ws4.Range("F33").Value = 845
$ws4.Range("F34").Value = 2552
$ws4.Range("F35").Value = 1028
$ws4.Range("F41").Value = 793
$ws4.Range("F43").Value = 1168
